$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data for "Olet" handling (basic OLET handling works only if target is a pipe)
$ws.Range("A4").Value = "Olet"
$ws.Range("A4").Font.Bold = $true

$ws.Range("B4").Value = "Pipe Types: Stålrør, sømløse, tap"

# Update the active selection to reflect where the user ended up editing
$ws.Range("B5").Select()
